$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the match data held in rows 47 and 48 (F,H,I,J,L,M,N,P,Q,R,T,U,V).
#    A,B,C,D,E,G,K,O,S stay as-is (they were already identical between
#    the two rows, or are untouched by the diff).
# ---------------------------------------------------------------------

$ws.Range("F47").Value = "Berkane"
$ws.Range("H47").Value = "Youssoufia Berrechid"
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 1.39
$ws.Range("L47").Value = 1.33
$ws.Range("M47").Value = "08/10/2023 18:54"
$ws.Range("N47").Value = 4.21
$ws.Range("P47").Value = 4.64
$ws.Range("Q47").Value = "08/10/2023 18:54"
$ws.Range("R47").Value = 7.84
$ws.Range("T47").Value = 9.800000000000001
$ws.Range("U47").Value = "08/10/2023 18:54"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/morocco/botola-pro/berkane-youssoufia-berrechid/W2sXokd9/"

$ws.Range("F48").Value = "Renaissance Zemamra"
$ws.Range("H48").Value = "Olympique de Safi"
$ws.Range("I48").Value = 2
$ws.Range("J48").Value = 2.98
$ws.Range("L48").Value = 3.95
$ws.Range("M48").Value = "08/10/2023 19:14"
$ws.Range("N48").Value = 2.85
$ws.Range("P48").Value = 2.75
$ws.Range("Q48").Value = "08/10/2023 19:14"
$ws.Range("R48").Value = 2.5
$ws.Range("T48").Value = 2.2
$ws.Range("U48").Value = "08/10/2023 19:14"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/morocco/botola-pro/renaissance-zemamra-olympique-de-safi/xSSTnVt3/"

# ---------------------------------------------------------------------
# 2) Append two new match rows (60, 61) after the existing last row (59).
#    Copy the cell formatting from the row above so the new rows carry
#    the same styles (bold/border on column A, datetime format on col E).
# ---------------------------------------------------------------------

function Add-MatchRow($RowIndex, $Indice, $Pais, $Torneio, $Temporada, $DataPartida, $Home, $HomeGols, $Away, $AwayGols, $HomeOpeningOdds, $HomeOpeningDataHora, $HomeClosingOdds, $HomeClosingDataHora, $DrawOpeningOdds, $DrawOpeningDataHora, $DrawClosingOdds, $DrawClosingDataHora, $AwayOpeningOdds, $AwayOpeningDataHora, $AwayClosingOdds, $AwayClosingDataHora, $UrlPartida) {
    $prevRow = $RowIndex - 1

    # Clone the number formatting/style from the row above first.
    $ws.Cells.Item($prevRow, 1).Copy()
    $ws.Cells.Item($RowIndex, 1).PasteSpecial(-4122)
    $ws.Cells.Item($prevRow, 5).Copy()
    $ws.Cells.Item($RowIndex, 5).PasteSpecial(-4122)

    $ws.Cells.Item($RowIndex, 1).Value = $Indice
    $ws.Cells.Item($RowIndex, 2).Value = $Pais
    $ws.Cells.Item($RowIndex, 3).Value = $Torneio
    $ws.Cells.Item($RowIndex, 4).Value = $Temporada
    $ws.Cells.Item($RowIndex, 5).Value = $DataPartida
    $ws.Cells.Item($RowIndex, 6).Value = $Home
    $ws.Cells.Item($RowIndex, 7).Value = $HomeGols
    $ws.Cells.Item($RowIndex, 8).Value = $Away
    $ws.Cells.Item($RowIndex, 9).Value = $AwayGols
    $ws.Cells.Item($RowIndex, 10).Value = $HomeOpeningOdds
    $ws.Cells.Item($RowIndex, 11).Value = $HomeOpeningDataHora
    $ws.Cells.Item($RowIndex, 12).Value = $HomeClosingOdds
    $ws.Cells.Item($RowIndex, 13).Value = $HomeClosingDataHora
    $ws.Cells.Item($RowIndex, 14).Value = $DrawOpeningOdds
    $ws.Cells.Item($RowIndex, 15).Value = $DrawOpeningDataHora
    $ws.Cells.Item($RowIndex, 16).Value = $DrawClosingOdds
    $ws.Cells.Item($RowIndex, 17).Value = $DrawClosingDataHora
    $ws.Cells.Item($RowIndex, 18).Value = $AwayOpeningOdds
    $ws.Cells.Item($RowIndex, 19).Value = $AwayOpeningDataHora
    $ws.Cells.Item($RowIndex, 20).Value = $AwayClosingOdds
    $ws.Cells.Item($RowIndex, 21).Value = $AwayClosingDataHora
    $ws.Cells.Item($RowIndex, 22).Value = $UrlPartida
}

Add-MatchRow 60 59 "morocco" "botola-pro" "2023-2024" `
    45235.66666666666 "IR Tanger" 2 "Maghreb Fez" 2 `
    2.71 "04/11/2023 04:12" `
    2.55 "05/11/2023 15:56" `
    3.03 "04/11/2023 04:12" `
    2.75 "05/11/2023 15:51" `
    2.47 "04/11/2023 04:12" `
    3.18 "05/11/2023 15:56" `
    "https://www.betexplorer.com/football/morocco/botola-pro/ir-tanger-maghreb-fez/trnEXILj/"

Add-MatchRow 61 60 "morocco" "botola-pro" "2023-2024" `
    45235.76041666666 "Mouloudia Oujda" 1 "Olympique de Safi" 1 `
    3.67 "04/11/2023 06:43" `
    3.75 "05/11/2023 11:23" `
    3.01 "04/11/2023 06:43" `
    2.89 "05/11/2023 16:18" `
    2 "04/11/2023 06:43" `
    2.16 "05/11/2023 11:23" `
    "https://www.betexplorer.com/football/morocco/botola-pro/mouloudia-oujda-olympique-de-safi/Sp7ld2bi/"
